$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("BD2").Value = 126

$ws.Range("Q5").Value = 2.1
$ws.Range("R5").Value = 1.7

$ws.Range("L6").Value = 9.5
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("AD6").Value = 9.5
$ws.Range("AG6").Value = 23
$ws.Range("AK6").Value = 81

$ws.Range("G14").Value = 5.5
$ws.Range("H14").Value = 3.7
$ws.Range("I14").Value = 1.53
$ws.Range("Q14").Value = 2
$ws.Range("R14").Value = 1.85
$ws.Range("W14").Value = 15
$ws.Range("AN14").Value = 7.5
$ws.Range("AW14").Value = 3.4

$ws.Range("G17").Value = 1.83
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 2.5
$ws.Range("N17").Value = 9.5
$ws.Range("X17").Value = 8.5
$ws.Range("Y17").Value = 8.5
$ws.Range("AG17").Value = 12
$ws.Range("AI17").Value = 15
$ws.Range("AN17").Value = 3.75
$ws.Range("AU17").Value = 8.5
$ws.Range("AX17").Value = 23
$ws.Range("BB17").Value = 251

$ws.Rows.Item(18).Delete()
